# dar de alta alumnos + solicitar correo con link de cambio de contraseña + cambio de contraseña
#
# Adds a new student ("Gallego Doncel, Aljenadro") as row 18 of the student
# list, with his DNI, a hyperlinked e-mail address (mailto: link, as Excel
# auto-creates when you type an e-mail address into a cell) and his
# convocatoria/matricula/movilidad/clases data - mirroring the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new row of data -------------------------------------------------
# Fill the student's name before the DNI so that new shared-strings are
# appended in the same order as the source workbook (name, DNI, email).
$ws.Range("A18").Value = 8
$ws.Range("C18").Value = "Gallego Doncel, Aljenadro"
$ws.Range("B18").Value = "80239121X"
$ws.Range("D18").Value = "UO285577@uniovi.es"
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = "No"
$ws.Range("H18").Value = "No"
$ws.Range("I18").Value = "Clases Expositivas-A"
$ws.Range("J18").Value = "Prácticas de Aula/Semin-02"
$ws.Range("K18").Value = "Prácticas de Laboratorio-01"
$ws.Range("L18").Value = "Tutorías Grupales-01"

# --- borders -----------------------------------------------------------
# Columns A:H pick up a left/right only border (same look the list had
# when the row was typed straight below the table).
foreach ($col in @("A","B","C","D","E","F","G","H")) {
    $cell = $ws.Range($col + "18")
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
}

# Columns I:L reuse the exact boxed-border look already used by the rows
# above (copy the format down from row 17 so the existing style is reused
# instead of a new, duplicate one being created).
$ws.Range("I17:L17").Copy()
$ws.Range("I18:L18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- e-mail hyperlink ----------------------------------------------------
# Excel auto-links an e-mail address typed into a cell with a mailto: link.
$ws.Hyperlinks.Add($ws.Range("D18"), "mailto:UO285577@uniovi.es")

# --- selection -----------------------------------------------------------
$ws.Range("D7").Select()
